# Append 23 new daily weather rows (366-388) to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 366
$ws.Range("F366").Value = 19
$ws.Range("G366").Value = 10
$ws.Range("H366").Value = 1
$ws.Range("I366").Value = -8
$ws.Range("J366").Value = -11
$ws.Range("K366").Value = -14
$ws.Range("L366").Value = 31
$ws.Range("M366").Value = 18
$ws.Range("N366").Value = 6
$ws.Range("O366").Value = 1016
$ws.Range("P366").Value = 1014
$ws.Range("Q366").Value = 1012
$ws.Range("R366").Value = 14
$ws.Range("S366").Value = 12
$ws.Range("T366").Value = 10
$ws.Range("U366").Value = 29
$ws.Range("V366").Value = 8
$ws.Range("X366").Value = 0
$ws.Range("Y366").Value = 4
$ws.Range("AA366").Value = 196

# Row 367
$ws.Range("F367").Value = 22
$ws.Range("G367").Value = 15
$ws.Range("H367").Value = 8
$ws.Range("I367").Value = 9
$ws.Range("J367").Value = -3
$ws.Range("K367").Value = -6
$ws.Range("L367").Value = 77
$ws.Range("M367").Value = 26
$ws.Range("N367").Value = 9
$ws.Range("O367").Value = 1016
$ws.Range("P367").Value = 1013
$ws.Range("Q367").Value = 1010
$ws.Range("R367").Value = 14
$ws.Range("S367").Value = 10
$ws.Range("T367").Value = 6
$ws.Range("U367").Value = 23
$ws.Range("V367").Value = 8
$ws.Range("X367").Value = 0
$ws.Range("Y367").Value = 2
$ws.Range("Z367").Value = "Rain"
$ws.Range("AA367").Value = 198

# Row 368
$ws.Range("F368").Value = 15
$ws.Range("G368").Value = 12
$ws.Range("H368").Value = 8
$ws.Range("I368").Value = 9
$ws.Range("J368").Value = 7
$ws.Range("K368").Value = 4
$ws.Range("L368").Value = 88
$ws.Range("M368").Value = 64
$ws.Range("N368").Value = 46
$ws.Range("O368").Value = 1014
$ws.Range("P368").Value = 1011
$ws.Range("Q368").Value = 1007
$ws.Range("R368").Value = 18
$ws.Range("S368").Value = 10
$ws.Range("T368").Value = 6
$ws.Range("U368").Value = 23
$ws.Range("V368").Value = 8
$ws.Range("X368").Value = 3.05
$ws.Range("Y368").Value = 6
$ws.Range("Z368").Value = "Rain-Snow"
$ws.Range("AA368").Value = 28

# Row 369
$ws.Range("F369").Value = 14
$ws.Range("G369").Value = 9
$ws.Range("H369").Value = 5
$ws.Range("I369").Value = 6
$ws.Range("J369").Value = 2
$ws.Range("K369").Value = -2
$ws.Range("L369").Value = 81
$ws.Range("M369").Value = 59
$ws.Range("N369").Value = 25
$ws.Range("O369").Value = 1015
$ws.Range("P369").Value = 1012
$ws.Range("Q369").Value = 1009
$ws.Range("R369").Value = 14
$ws.Range("S369").Value = 12
$ws.Range("T369").Value = 10
$ws.Range("U369").Value = 23
$ws.Range("V369").Value = 5
$ws.Range("X369").Value = 0
$ws.Range("Y369").Value = 3
$ws.Range("AA369").Value = 37

# Row 370
$ws.Range("F370").Value = 18
$ws.Range("G370").Value = 11
$ws.Range("H370").Value = 4
$ws.Range("I370").Value = 2
$ws.Range("J370").Value = -1
$ws.Range("K370").Value = -6
$ws.Range("L370").Value = 71
$ws.Range("M370").Value = 39
$ws.Range("N370").Value = 10
$ws.Range("O370").Value = 1021
$ws.Range("P370").Value = 1019
$ws.Range("Q370").Value = 1014
$ws.Range("R370").Value = 19
$ws.Range("S370").Value = 13
$ws.Range("T370").Value = 10
$ws.Range("U370").Value = 14
$ws.Range("V370").Value = 5
$ws.Range("X370").Value = 0
$ws.Range("Y370").Value = 5
$ws.Range("AA370").Value = 243

# Row 371
$ws.Range("F371").Value = 20
$ws.Range("G371").Value = 14
$ws.Range("H371").Value = 8
$ws.Range("I371").Value = 0
$ws.Range("J371").Value = -3
$ws.Range("K371").Value = -6
$ws.Range("L371").Value = 47
$ws.Range("M371").Value = 28
$ws.Range("N371").Value = 10
$ws.Range("O371").Value = 1022
$ws.Range("P371").Value = 1019
$ws.Range("Q371").Value = 1014
$ws.Range("R371").Value = 14
$ws.Range("S371").Value = 13
$ws.Range("T371").Value = 10
$ws.Range("U371").Value = 8
$ws.Range("V371").Value = 2
$ws.Range("X371").Value = 0
$ws.Range("Y371").Value = 4
$ws.Range("AA371").Value = 12

# Row 372
$ws.Range("F372").Value = 19
$ws.Range("G372").Value = 16
$ws.Range("H372").Value = 12
$ws.Range("I372").Value = 8
$ws.Range("J372").Value = 1
$ws.Range("K372").Value = -2
$ws.Range("L372").Value = 67
$ws.Range("M372").Value = 34
$ws.Range("N372").Value = 15
$ws.Range("O372").Value = 1018
$ws.Range("P372").Value = 1016
$ws.Range("Q372").Value = 1014
$ws.Range("R372").Value = 10
$ws.Range("S372").Value = 9
$ws.Range("T372").Value = 3
$ws.Range("U372").Value = 14
$ws.Range("V372").Value = 5
$ws.Range("X372").Value = 0.25
$ws.Range("Y372").Value = 8
$ws.Range("Z372").Value = "Rain"
$ws.Range("AA372").Value = 258

# Row 373
$ws.Range("F373").Value = 19
$ws.Range("G373").Value = 15
$ws.Range("H373").Value = 11
$ws.Range("I373").Value = 9
$ws.Range("J373").Value = 7
$ws.Range("K373").Value = 6
$ws.Range("L373").Value = 77
$ws.Range("M373").Value = 57
$ws.Range("N373").Value = 35
$ws.Range("O373").Value = 1018
$ws.Range("P373").Value = 1012
$ws.Range("Q373").Value = 1007
$ws.Range("R373").Value = 11
$ws.Range("S373").Value = 10
$ws.Range("T373").Value = 8
$ws.Range("U373").Value = 23
$ws.Range("V373").Value = 6
$ws.Range("X373").Value = 0.51
$ws.Range("Y373").Value = 8
$ws.Range("Z373").Value = "Rain"
$ws.Range("AA373").Value = 283

# Row 374
$ws.Range("F374").Value = 23
$ws.Range("G374").Value = 17
$ws.Range("H374").Value = 9
$ws.Range("I374").Value = 8
$ws.Range("J374").Value = 5
$ws.Range("K374").Value = 0
$ws.Range("L374").Value = 82
$ws.Range("M374").Value = 45
$ws.Range("N374").Value = 12
$ws.Range("O374").Value = 1012
$ws.Range("P374").Value = 1009
$ws.Range("Q374").Value = 1004
$ws.Range("R374").Value = 18
$ws.Range("S374").Value = 11
$ws.Range("T374").Value = 10
$ws.Range("U374").Value = 26
$ws.Range("V374").Value = 8
$ws.Range("X374").Value = 0
$ws.Range("Y374").Value = 4
$ws.Range("Z374").Value = "Rain"
$ws.Range("AA374").Value = 184

# Row 375
$ws.Range("F375").Value = 24
$ws.Range("G375").Value = 19
$ws.Range("H375").Value = 13
$ws.Range("I375").Value = 9
$ws.Range("J375").Value = 4
$ws.Range("K375").Value = -3
$ws.Range("L375").Value = 72
$ws.Range("M375").Value = 37
$ws.Range("N375").Value = 11
$ws.Range("O375").Value = 1009
$ws.Range("P375").Value = 1006
$ws.Range("Q375").Value = 1001
$ws.Range("R375").Value = 11
$ws.Range("S375").Value = 10
$ws.Range("T375").Value = 10
$ws.Range("U375").Value = 26
$ws.Range("V375").Value = 8
$ws.Range("X375").Value = 0
$ws.Range("Y375").Value = 4
$ws.Range("Z375").Value = "Rain"
$ws.Range("AA375").Value = 310

# Row 376
$ws.Range("F376").Value = 19
$ws.Range("G376").Value = 14
$ws.Range("H376").Value = 8
$ws.Range("I376").Value = 0
$ws.Range("J376").Value = -7
$ws.Range("K376").Value = -20
$ws.Range("L376").Value = 41
$ws.Range("M376").Value = 23
$ws.Range("N376").Value = 4
$ws.Range("O376").Value = 1014
$ws.Range("P376").Value = 1009
$ws.Range("Q376").Value = 1004
$ws.Range("R376").Value = 14
$ws.Range("S376").Value = 12
$ws.Range("T376").Value = 10
$ws.Range("U376").Value = 35
$ws.Range("V376").Value = 11
$ws.Range("X376").Value = 0
$ws.Range("Y376").Value = 2
$ws.Range("AA376").Value = 309

# Row 377
$ws.Range("F377").Value = 18
$ws.Range("G377").Value = 12
$ws.Range("H377").Value = 5
$ws.Range("I377").Value = -4
$ws.Range("J377").Value = -7
$ws.Range("K377").Value = -12
$ws.Range("L377").Value = 43
$ws.Range("M377").Value = 23
$ws.Range("N377").Value = 6
$ws.Range("O377").Value = 1018
$ws.Range("P377").Value = 1015
$ws.Range("Q377").Value = 1011
$ws.Range("R377").Value = 19
$ws.Range("S377").Value = 15
$ws.Range("T377").Value = 10
$ws.Range("U377").Value = 18
$ws.Range("V377").Value = 5
$ws.Range("X377").Value = 0
$ws.Range("Y377").Value = 2
$ws.Range("AA377").Value = 360

# Row 378
$ws.Range("F378").Value = 12
$ws.Range("G378").Value = 11
$ws.Range("H378").Value = 9
$ws.Range("I378").Value = -2
$ws.Range("J378").Value = -6
$ws.Range("K378").Value = -9
$ws.Range("L378").Value = 41
$ws.Range("M378").Value = 26
$ws.Range("N378").Value = 14
$ws.Range("O378").Value = 1019
$ws.Range("P378").Value = 1015
$ws.Range("Q378").Value = 1010
$ws.Range("R378").Value = 14
$ws.Range("S378").Value = 12
$ws.Range("T378").Value = 10
$ws.Range("U378").Value = 35
$ws.Range("V378").Value = 19
$ws.Range("X378").Value = 0
$ws.Range("Y378").Value = 8
$ws.Range("AA378").Value = 357

# Row 379
$ws.Range("F379").Value = 15
$ws.Range("G379").Value = 10
$ws.Range("H379").Value = 5
$ws.Range("I379").Value = 0
$ws.Range("J379").Value = -4
$ws.Range("K379").Value = -9
$ws.Range("L379").Value = 57
$ws.Range("M379").Value = 32
$ws.Range("N379").Value = 18
$ws.Range("O379").Value = 1021
$ws.Range("P379").Value = 1018
$ws.Range("Q379").Value = 1016
$ws.Range("R379").Value = 14
$ws.Range("S379").Value = 12
$ws.Range("T379").Value = 10
$ws.Range("U379").Value = 19
$ws.Range("V379").Value = 6
$ws.Range("X379").Value = 0
$ws.Range("Y379").Value = 3
$ws.Range("AA379").Value = 172

# Row 380
$ws.Range("F380").Value = 19
$ws.Range("G380").Value = 12
$ws.Range("H380").Value = 5
$ws.Range("I380").Value = -4
$ws.Range("J380").Value = -7
$ws.Range("K380").Value = -10
$ws.Range("L380").Value = 46
$ws.Range("M380").Value = 25
$ws.Range("N380").Value = 7
$ws.Range("O380").Value = 1018
$ws.Range("P380").Value = 1016
$ws.Range("Q380").Value = 1013
$ws.Range("R380").Value = 14
$ws.Range("S380").Value = 12
$ws.Range("T380").Value = 10
$ws.Range("U380").Value = 19
$ws.Range("V380").Value = 5
$ws.Range("X380").Value = 0
$ws.Range("Y380").Value = 6
$ws.Range("AA380").Value = 154

# Row 381
$ws.Range("F381").Value = 21
$ws.Range("G381").Value = 16
$ws.Range("H381").Value = 9
$ws.Range("I381").Value = -4
$ws.Range("J381").Value = -6
$ws.Range("K381").Value = -8
$ws.Range("L381").Value = 30
$ws.Range("M381").Value = 19
$ws.Range("N381").Value = 7
$ws.Range("O381").Value = 1016
$ws.Range("P381").Value = 1013
$ws.Range("Q381").Value = 1009
$ws.Range("R381").Value = 14
$ws.Range("S381").Value = 12
$ws.Range("T381").Value = 10
$ws.Range("U381").Value = 26
$ws.Range("V381").Value = 8
$ws.Range("X381").Value = 0
$ws.Range("Y381").Value = 2
$ws.Range("AA381").Value = 293

# Row 382
$ws.Range("F382").Value = 21
$ws.Range("G382").Value = 14
$ws.Range("H382").Value = 8
$ws.Range("I382").Value = -2
$ws.Range("J382").Value = -4
$ws.Range("K382").Value = -6
$ws.Range("L382").Value = 37
$ws.Range("M382").Value = 25
$ws.Range("N382").Value = 10
$ws.Range("O382").Value = 1016
$ws.Range("P382").Value = 1013
$ws.Range("Q382").Value = 1010
$ws.Range("R382").Value = 14
$ws.Range("S382").Value = 12
$ws.Range("T382").Value = 10
$ws.Range("U382").Value = 11
$ws.Range("V382").Value = 3
$ws.Range("X382").Value = 0
$ws.Range("Y382").Value = 3
$ws.Range("AA382").Value = 116

# Row 383
$ws.Range("F383").Value = 24
$ws.Range("G383").Value = 17
$ws.Range("H383").Value = 9
$ws.Range("I383").Value = -3
$ws.Range("J383").Value = -5
$ws.Range("K383").Value = -8
$ws.Range("L383").Value = 37
$ws.Range("M383").Value = 19
$ws.Range("N383").Value = 6
$ws.Range("O383").Value = 1016
$ws.Range("P383").Value = 1013
$ws.Range("Q383").Value = 1009
$ws.Range("R383").Value = 18
$ws.Range("S383").Value = 12
$ws.Range("T383").Value = 10
$ws.Range("U383").Value = 19
$ws.Range("V383").Value = 5
$ws.Range("X383").Value = 0
$ws.Range("Y383").Value = 3
$ws.Range("AA383").Value = 216

# Row 384
$ws.Range("F384").Value = 25
$ws.Range("G384").Value = 18
$ws.Range("H384").Value = 12
$ws.Range("I384").Value = -3
$ws.Range("J384").Value = -5
$ws.Range("K384").Value = -7
$ws.Range("L384").Value = 31
$ws.Range("M384").Value = 16
$ws.Range("N384").Value = 7
$ws.Range("O384").Value = 1016
$ws.Range("P384").Value = 1013
$ws.Range("Q384").Value = 1009
$ws.Range("R384").Value = 14
$ws.Range("S384").Value = 12
$ws.Range("T384").Value = 10
$ws.Range("U384").Value = 19
$ws.Range("V384").Value = 6
$ws.Range("X384").Value = 0
$ws.Range("Y384").Value = 2
$ws.Range("AA384").Value = 299

# Row 385
$ws.Range("F385").Value = 23
$ws.Range("G385").Value = 18
$ws.Range("H385").Value = 12
$ws.Range("I385").Value = 2
$ws.Range("J385").Value = -3
$ws.Range("K385").Value = -6
$ws.Range("L385").Value = 34
$ws.Range("M385").Value = 21
$ws.Range("N385").Value = 10
$ws.Range("O385").Value = 1018
$ws.Range("P385").Value = 1015
$ws.Range("Q385").Value = 1010
$ws.Range("R385").Value = 14
$ws.Range("S385").Value = 11
$ws.Range("T385").Value = 10
$ws.Range("U385").Value = 35
$ws.Range("V385").Value = 6
$ws.Range("X385").Value = 0
$ws.Range("Y385").Value = 3
$ws.Range("AA385").Value = 304

# Row 386
$ws.Range("F386").Value = 19
$ws.Range("G386").Value = 16
$ws.Range("H386").Value = 12
$ws.Range("I386").Value = 8
$ws.Range("J386").Value = 7
$ws.Range("K386").Value = 3
$ws.Range("L386").Value = 72
$ws.Range("M386").Value = 54
$ws.Range("N386").Value = 33
$ws.Range("O386").Value = 1021
$ws.Range("P386").Value = 1017
$ws.Range("Q386").Value = 1014
$ws.Range("R386").Value = 11
$ws.Range("S386").Value = 10
$ws.Range("T386").Value = 7
$ws.Range("U386").Value = 23
$ws.Range("V386").Value = 8
$ws.Range("X386").Value = 3.05
$ws.Range("Y386").Value = 3
$ws.Range("Z386").Value = "Rain-Thunderstorm"
$ws.Range("AA386").Value = 61

# Row 387
$ws.Range("F387").Value = 23
$ws.Range("G387").Value = 17
$ws.Range("H387").Value = 9
$ws.Range("I387").Value = 7
$ws.Range("J387").Value = 3
$ws.Range("K387").Value = 0
$ws.Range("L387").Value = 67
$ws.Range("M387").Value = 38
$ws.Range("N387").Value = 12
$ws.Range("O387").Value = 1018
$ws.Range("P387").Value = 1016
$ws.Range("Q387").Value = 1011
$ws.Range("R387").Value = 14
$ws.Range("S387").Value = 11
$ws.Range("T387").Value = 10
$ws.Range("U387").Value = 11
$ws.Range("V387").Value = 5
$ws.Range("X387").Value = 0
$ws.Range("Y387").Value = 3
$ws.Range("AA387").Value = 284

# Row 388
$ws.Range("F388").Value = 24
$ws.Range("G388").Value = 19
$ws.Range("H388").Value = 14
$ws.Range("I388").Value = 5
$ws.Range("J388").Value = 3
$ws.Range("K388").Value = 2
$ws.Range("L388").Value = 48
$ws.Range("M388").Value = 31
$ws.Range("N388").Value = 14
$ws.Range("O388").Value = 1016
$ws.Range("P388").Value = 1013
$ws.Range("Q388").Value = 1008
$ws.Range("R388").Value = 18
$ws.Range("S388").Value = 11
$ws.Range("T388").Value = 10
$ws.Range("U388").Value = 14
$ws.Range("V388").Value = 3
$ws.Range("X388").Value = 0
$ws.Range("Y388").Value = 5
$ws.Range("AA388").Value = 175

# Restore a normal (unscrolled) view with a single active cell, as in the saved workbook.
$ws.Range("O5").Select()
